$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 10 (pushes former rows 10-18 down to 11-19) and
# populate it with the new "8b" task for Ilkay.
$ws.Rows(10).Insert()

$ws.Range("A10").Value = "8b"
$ws.Range("B10").Value = "Create user stories and add it to doc"
$ws.Range("C10").Value = "Ilkay"
$ws.Range("D10").Value = "Task 4 is finished"
$ws.Range("E10").Value = "1. team review 2. uploaded to github 3. present to a team"
$ws.Range("F10").Value = "New"
$ws.Range("F10").WrapText = $true
$ws.Rows(10).RowHeight = 30

# Add the new "Time" column (G) with estimates for the first few tasks.
$ws.Range("G2").Value = "1 day"
$ws.Range("G2").WrapText = $true

$ws.Range("G3").Value = "7  days"
$ws.Range("G3").WrapText = $true

$ws.Range("G4").Value = "1 hour"
$ws.Range("G4").WrapText = $true

$ws.Range("G1").Value = "Time"
$ws.Range("G1").WrapText = $true
$ws.Range("G1").Font.Bold = $true

[void]$ws.Range("G5").Select()
